# Insert a new data row at row 243 (pushing existing rows 243..282 down to
# 244..283) and populate it with the new record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new entire row above the current row 243.
$ws.Rows.Item(243).EntireRow.Insert()

# Populate the newly inserted row 243 with the new record's data.
$ws.Cells.Item(243, 1).Value = 9
$ws.Cells.Item(243, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(243, 3).Value = "Metropolitana"
$ws.Cells.Item(243, 4).Value = 44504
$ws.Cells.Item(243, 5).Value = 13
$ws.Cells.Item(243, 6).Value = 100112031
$ws.Cells.Item(243, 7).Value = "Poroto verde"
$ws.Cells.Item(243, 8).Value = "Magnum"
$ws.Cells.Item(243, 9).Value = "Primera"
$ws.Cells.Item(243, 10).Value = 70
$ws.Cells.Item(243, 11).Value = 38000
$ws.Cells.Item(243, 12).Value = 41000
$ws.Cells.Item(243, 13).Value = 39500
$ws.Cells.Item(243, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(243, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(243, 16).Value = 1580
$ws.Cells.Item(243, 17).Value = 25
$ws.Cells.Item(243, 18).Value = "Hortaliza"
